$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.783.96"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "3.437.79"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.695"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("D15").Value = "3.446.28"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "62.674.69"
$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000155"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +19.83%  "

$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "315.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("E25").Value = "  +1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("E28").Value = "  +5.82%  "

$ws.Range("E29").Value = "  +9.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.93%  "

$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0485"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "143.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.55%  "

$ws.Range("E42").Value = "  +1.47%  "

$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "

$ws.Range("D48").Value = "2.110.15"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("B49").Value = "OceanProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +32.17%  "

$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.68%  "
